$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# K2:K19 hold the librarySampleNumber "E7760" -> rename to "E7420".
# Every row shares the same text, so update each cell individually
# (keeps the shared-string table / formula encoding as tidy as possible).
for ($r = 2; $r -le 19; $r++) {
    $ws.Cells.Item($r, 11).Value = "E7420"
}

# L2:L19 were a literal boolean FALSE; turn them into an explicit
# =FALSE() formula (still evaluating to FALSE/0), one cell at a time so
# Excel doesn't collapse them into a single shared-formula group.
for ($r = 2; $r -le 19; $r++) {
    $ws.Cells.Item($r, 12).Formula = "=FALSE()"
}

# Selection moves from L2:L19 to K2:K19.
$ws.Range("K2:K19").Select()
